$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.610.73'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.457.09'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '564.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('E9').Value = '  -5.30%  '
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.344'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.04%  '
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '2.911.26'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '68.594.55'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.42%  '
$ws.Range('D17').Value = '2.465.49'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '353.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.37%  '
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('E23').Value = '  -4.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('E25').Value = '  -4.59%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.578.97'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').Value = '  -5.92%  '
$ws.Range('D29').Value = '0.0₃0841'
$ws.Range('E29').Value = '  -5.37%  '
$ws.Range('E30').Value = '  -6.62%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '431.32'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.35%  '
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +105.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('E38').Value = '  -5.49%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.85%  '
$ws.Range('E47').Value = '  -2.72%  '
$ws.Range('E48').Value = '  -5.14%  '
$ws.Range('E49').Value = '  -2.08%  '
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0917'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.08%  '
